$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.36880911882937539,
    -0.070437720477297461,
    -0.0089999987454785213,
    0.061990161506578545,
    -0.0059999987647394448,
    -0.0059999987280932032,
    -0.019999998528499319,
    -0.019999998524282248,
    -0.0059999987205934247,
    -0.0059999987214141015,
    0.0086165535474442834,
    -0.0059999987212959738,
    -0.0059999987240253461,
    -0.011999998641210041,
    -0.0059999987295231705,
    -0.0059999987356857964,
    -0.005999998749890878,
    0.044351738860677337,
    -0.086494401496521434,
    -0.0089999987542199733,
    -0.0089999987525102298,
    -0.0089999987515083646,
    -0.0089999987355229294,
    -0.032735585952955226,
    -0.041999998242409298,
    -0.0059999987232046692,
    -0.0059999987164998103,
    -0.0059999986885896917,
    -0.011999998585226379,
    -0.019999998463712476,
    -0.014999998521512126,
    -0.020999998435083711,
    -0.0059999986464260857
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
